# Generate Report for Handoff
# Adds two new tracked files (dependency PNGs) around the existing markdown
# hand-off row on each of the three report sheets (Overview, zh-cn, de-de),
# refreshes the existing row's values/timestamps, and rewires the hyperlinks.

$wb = $excel.ActiveWorkbook

# Cornflower blue (FF6495ED) expressed as a BGR OLE_COLOR for Font.Color.
$hyperlinkColor = 15570276
$dateFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "0dea70fd-23e2-4c6b-8ca4-3324daab77fd.png"
$ov.Range("B2").Value = "Ready for handoff"
$ov.Range("C2").Value = "Ready for handoff"
$ov.Range("D2").Value = "2016-03-23 11:09:25"

$ov.Range("A3").Value = "3322df89-3800-43de-ba07-a91ca2f3c5ca.md"
$ov.Range("B3").Value = "Ready for handoff"
$ov.Range("C3").Value = "Ready for handoff"
$ov.Range("D3").Value = "2016-03-23 11:09:25"

$ov.Range("A4").Value = "c4e015c0-c120-4d16-ae74-3ef426950a9d.png"
$ov.Range("B4").Value = "Ready for handoff"
$ov.Range("C4").Value = "Ready for handoff"
$ov.Range("D4").Value = "2016-03-23 11:09:25"

$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/ae4718d54c9fa7bb9b78747a87dfe0dbd24f098a/e2e/0dea70fd-23e2-4c6b-8ca4-3324daab77fd.png", "", "", "0dea70fd-23e2-4c6b-8ca4-3324daab77fd.png")
$ov.Hyperlinks.Add($ov.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/ae4718d54c9fa7bb9b78747a87dfe0dbd24f098a/e2e/3322df89-3800-43de-ba07-a91ca2f3c5ca.md", "", "", "3322df89-3800-43de-ba07-a91ca2f3c5ca.md")
$ov.Hyperlinks.Add($ov.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/ae4718d54c9fa7bb9b78747a87dfe0dbd24f098a/e2e/c4e015c0-c120-4d16-ae74-3ef426950a9d.png", "", "", "c4e015c0-c120-4d16-ae74-3ef426950a9d.png")

$ov.Range("A2:A4").Font.Underline = 2
$ov.Range("A2:A4").Font.Color = $hyperlinkColor
$ov.Range("D2:D4").NumberFormat = $dateFormat

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "0dea70fd-23e2-4c6b-8ca4-3324daab77fd.png"
$zh.Range("B2").Value = ".png"
$zh.Range("C2").Value = "Ready for handoff"
$zh.Range("D2").Value = "5bba332a991b20e9cd3efda468189a1f8e987b42.png"
$zh.Range("E2").Value = "2016-03-23 11:09:20"
$zh.Range("H2").Value = "0001-01-01 00:00:00"
$zh.Range("J2").Value = "IsDependency"
$zh.Range("K2").Value = "e2e\3322df89-3800-43de-ba07-a91ca2f3c5ca.md"

$zh.Range("A3").Value = "3322df89-3800-43de-ba07-a91ca2f3c5ca.md"
$zh.Range("B3").Value = ".md"
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("D3").Value = "3322df89-3800-43de-ba07-a91ca2f3c5ca.6831b018e6c9182276a120ebc19904d6eb316969.zh-cn.xlf"
$zh.Range("E3").Value = "2016-03-23 11:09:20"
$zh.Range("H3").Value = "0001-01-01 00:00:00"
$zh.Range("J3").Value = "Include"

$zh.Range("A4").Value = "c4e015c0-c120-4d16-ae74-3ef426950a9d.png"
$zh.Range("B4").Value = ".png"
$zh.Range("C4").Value = "Ready for handoff"
$zh.Range("D4").Value = "62491b4f22d0c51760fa0581ef68d2e645ec440e.png"
$zh.Range("E4").Value = "2016-03-23 11:09:20"
$zh.Range("H4").Value = "0001-01-01 00:00:00"
$zh.Range("J4").Value = "IsDependency"
$zh.Range("K4").Value = "e2e\3322df89-3800-43de-ba07-a91ca2f3c5ca.md"

$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/ae4718d54c9fa7bb9b78747a87dfe0dbd24f098a/e2e/0dea70fd-23e2-4c6b-8ca4-3324daab77fd.png", "", "", "0dea70fd-23e2-4c6b-8ca4-3324daab77fd.png")
$zh.Hyperlinks.Add($zh.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/82b3b4bfa663deba68eb7b63d8a5d5d2397f92b3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/5bba332a991b20e9cd3efda468189a1f8e987b42.png", "", "", "5bba332a991b20e9cd3efda468189a1f8e987b42.png")
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/ae4718d54c9fa7bb9b78747a87dfe0dbd24f098a/e2e/3322df89-3800-43de-ba07-a91ca2f3c5ca.md", "", "", "3322df89-3800-43de-ba07-a91ca2f3c5ca.md")
$zh.Hyperlinks.Add($zh.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/82b3b4bfa663deba68eb7b63d8a5d5d2397f92b3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/3322df89-3800-43de-ba07-a91ca2f3c5ca.6831b018e6c9182276a120ebc19904d6eb316969.zh-cn.xlf", "", "", "3322df89-3800-43de-ba07-a91ca2f3c5ca.6831b018e6c9182276a120ebc19904d6eb316969.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/ae4718d54c9fa7bb9b78747a87dfe0dbd24f098a/e2e/c4e015c0-c120-4d16-ae74-3ef426950a9d.png", "", "", "c4e015c0-c120-4d16-ae74-3ef426950a9d.png")
$zh.Hyperlinks.Add($zh.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/82b3b4bfa663deba68eb7b63d8a5d5d2397f92b3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/62491b4f22d0c51760fa0581ef68d2e645ec440e.png", "", "", "62491b4f22d0c51760fa0581ef68d2e645ec440e.png")

$zh.Range("A2:A4").Font.Underline = 2
$zh.Range("A2:A4").Font.Color = $hyperlinkColor
$zh.Range("D2:D4").Font.Underline = 2
$zh.Range("D2:D4").Font.Color = $hyperlinkColor
$zh.Range("E2:E4").NumberFormat = $dateFormat
$zh.Range("H2:H4").NumberFormat = $dateFormat

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "0dea70fd-23e2-4c6b-8ca4-3324daab77fd.png"
$de.Range("B2").Value = ".png"
$de.Range("C2").Value = "Ready for handoff"
$de.Range("D2").Value = "5bba332a991b20e9cd3efda468189a1f8e987b42.png"
$de.Range("E2").Value = "2016-03-23 11:09:25"
$de.Range("H2").Value = "0001-01-01 00:00:00"
$de.Range("J2").Value = "IsDependency"
$de.Range("K2").Value = "e2e\3322df89-3800-43de-ba07-a91ca2f3c5ca.md"

$de.Range("A3").Value = "3322df89-3800-43de-ba07-a91ca2f3c5ca.md"
$de.Range("B3").Value = ".md"
$de.Range("C3").Value = "Ready for handoff"
$de.Range("D3").Value = "3322df89-3800-43de-ba07-a91ca2f3c5ca.6831b018e6c9182276a120ebc19904d6eb316969.de-de.xlf"
$de.Range("E3").Value = "2016-03-23 11:09:25"
$de.Range("H3").Value = "0001-01-01 00:00:00"
$de.Range("J3").Value = "Include"

$de.Range("A4").Value = "c4e015c0-c120-4d16-ae74-3ef426950a9d.png"
$de.Range("B4").Value = ".png"
$de.Range("C4").Value = "Ready for handoff"
$de.Range("D4").Value = "62491b4f22d0c51760fa0581ef68d2e645ec440e.png"
$de.Range("E4").Value = "2016-03-23 11:09:25"
$de.Range("H4").Value = "0001-01-01 00:00:00"
$de.Range("J4").Value = "IsDependency"
$de.Range("K4").Value = "e2e\3322df89-3800-43de-ba07-a91ca2f3c5ca.md"

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/ae4718d54c9fa7bb9b78747a87dfe0dbd24f098a/e2e/0dea70fd-23e2-4c6b-8ca4-3324daab77fd.png", "", "", "0dea70fd-23e2-4c6b-8ca4-3324daab77fd.png")
$de.Hyperlinks.Add($de.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3d19121d8b665f5aa37ad88daed63d4d8787d3a7/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/5bba332a991b20e9cd3efda468189a1f8e987b42.png", "", "", "5bba332a991b20e9cd3efda468189a1f8e987b42.png")
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/ae4718d54c9fa7bb9b78747a87dfe0dbd24f098a/e2e/3322df89-3800-43de-ba07-a91ca2f3c5ca.md", "", "", "3322df89-3800-43de-ba07-a91ca2f3c5ca.md")
$de.Hyperlinks.Add($de.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3d19121d8b665f5aa37ad88daed63d4d8787d3a7/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/3322df89-3800-43de-ba07-a91ca2f3c5ca.6831b018e6c9182276a120ebc19904d6eb316969.de-de.xlf", "", "", "3322df89-3800-43de-ba07-a91ca2f3c5ca.6831b018e6c9182276a120ebc19904d6eb316969.de-de.xlf")
$de.Hyperlinks.Add($de.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/ae4718d54c9fa7bb9b78747a87dfe0dbd24f098a/e2e/c4e015c0-c120-4d16-ae74-3ef426950a9d.png", "", "", "c4e015c0-c120-4d16-ae74-3ef426950a9d.png")
$de.Hyperlinks.Add($de.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3d19121d8b665f5aa37ad88daed63d4d8787d3a7/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/62491b4f22d0c51760fa0581ef68d2e645ec440e.png", "", "", "62491b4f22d0c51760fa0581ef68d2e645ec440e.png")

$de.Range("A2:A4").Font.Underline = 2
$de.Range("A2:A4").Font.Color = $hyperlinkColor
$de.Range("D2:D4").Font.Underline = 2
$de.Range("D2:D4").Font.Color = $hyperlinkColor
$de.Range("E2:E4").NumberFormat = $dateFormat
$de.Range("H2:H4").NumberFormat = $dateFormat
